$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.811.71"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "3.315.56"
$ws.Range("E3").Value = "  +6.47%  "
$ws.Range("E4").Value = "  -0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "601.75"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.98"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.23%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.314.65"
$ws.Range("E8").Value = "  +6.70%  "
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("E10").Value = "  +3.32%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.57"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +5.76%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.474"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +4.40%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000249"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.50%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "34.85"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").Value = "3.861.15"
$ws.Range("E15").Value = "  +6.45%  "
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "3.315.06"
$ws.Range("E17").Value = "  +6.31%  "
$ws.Range("D18").Value = "63.921.79"
$ws.Range("E18").Value = "  +1.52%  "
$ws.Range("E19").Value = "  +3.96%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "481.87"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.14%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.24"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.47%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.735"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +5.89%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "8.05"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +4.99%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "13.53"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +4.93%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "84.88"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +1.87%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.28"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +5.42%  "
$ws.Range("E29").Value = "  -0.15%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "8.17"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +3.66%  "
$ws.Range("E31").Value = "  +4.54%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "29.43"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +10.59%  "
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E34").Value = "  +1.33%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.10"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.46%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.99"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.55%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "52.97"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.84%  "
$ws.Range("D38").Value = "0.0₃0755"
$ws.Range("E38").Value = "  +7.91%  "
$ws.Range("E39").Value = "  +4.47%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "432.34"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.65%  "
$ws.Range("D41").Value = "3.051.09"
$ws.Range("E41").Value = "  +5.29%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "8.42"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.10%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.76"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.92%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.113"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("E46").Value = "  +5.01%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "26.47"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.10%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "36.00"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +14.63%  "
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("E51").Value = "  +3.26%  "
